$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 335067
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 335067
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 1005201
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -1011335
$ws.Range("H113").Value = 7000.8335
$ws.Range("I113").Value = 5001.6665
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 5001.6665
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -1747.6665
$ws.Range("N113").Value = -15508
$ws.Range("H132").Value = 2140.6128
$ws.Range("I132").Value = 1288.5769
$ws.Range("K132").Value = 3865.7307
$ws.Range("M132").Value = -1335.7307
$ws.Range("H138").Value = 5173.763
$ws.Range("I138").Value = 968.6429000000001
$ws.Range("J138").Value = 7626.75
$ws.Range("K138").Value = 2905.9287
$ws.Range("L138").Value = 22880.25
$ws.Range("M138").Value = 2234.0713
$ws.Range("N138").Value = -33160.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5808.6924
$ws.Range("I61").Value = 6782.2383
$ws.Range("K61").Value = 6782.2383
$ws.Range("M61").Value = -6570.2383
$ws.Range("H102").Value = 7410567
$ws.Range("H109").Value = 46000
$ws.Range("J109").Value = 46000
$ws.Range("L109").Value = 46000
$ws.Range("N109").Value = -48774
$ws.Range("H110").Value = 1033.3334
$ws.Range("I110").Value = 1700
$ws.Range("K110").Value = 1700
$ws.Range("M110").Value = 345
$ws.Range("H136").Value = 5808.6924
$ws.Range("I136").Value = 6782.2383
$ws.Range("K136").Value = 20346.7149
$ws.Range("M136").Value = -17796.7149

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 16668314
$ws.Range("I86").Value = 25642576
$ws.Range("J86").Value = 1828.5714
$ws.Range("K86").Value = 25642576
$ws.Range("L86").Value = 1828.5714
$ws.Range("M86").Value = -25641453
$ws.Range("N86").Value = -4074.5714
$ws.Range("H89").Value = 16668314
$ws.Range("I89").Value = 25642576
$ws.Range("J89").Value = 1828.5714
$ws.Range("K89").Value = 128212880
$ws.Range("L89").Value = 9142.857
$ws.Range("M89").Value = -128207264
$ws.Range("N89").Value = -20374.857
$ws.Range("H94").Value = 1995.9231
$ws.Range("I94").Value = 1759
$ws.Range("K94").Value = 1759
$ws.Range("M94").Value = -1308
$ws.Range("H99").Value = 500001250
$ws.Range("I99").Value = 1000000000
$ws.Range("K99").Value = 1000000000
$ws.Range("M99").Value = -999998502
$ws.Range("H105").Value = 24872.666
$ws.Range("I105").Value = 35301.668
$ws.Range("J105").Value = 4014.6667
$ws.Range("K105").Value = 35301.668
$ws.Range("L105").Value = 4014.6667
$ws.Range("M105").Value = -33554.668
$ws.Range("N105").Value = -7508.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 23500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 23500
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -24972
$ws.Range("H61").Value = 23500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 23500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 23500
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -24196
$ws.Range("H62").Value = 6121.75
$ws.Range("I62").Value = 8661.666999999999
$ws.Range("J62").Value = 4597.8
$ws.Range("K62").Value = 8661.666999999999
$ws.Range("L62").Value = 4597.8
$ws.Range("M62").Value = -8037.666999999999
$ws.Range("N62").Value = -5845.8
$ws.Range("H65").Value = 6121.75
$ws.Range("I65").Value = 8661.666999999999
$ws.Range("J65").Value = 4597.8
$ws.Range("K65").Value = 43308.335
$ws.Range("L65").Value = 22989
$ws.Range("M65").Value = -40188.335
$ws.Range("N65").Value = -29229
$ws.Range("H99").Value = 3639.5386
$ws.Range("J99").Value = 7362.8
$ws.Range("L99").Value = 7362.8
$ws.Range("N99").Value = -10358.8
$ws.Range("H105").Value = 10000
$ws.Range("I105").Value = 10000
$ws.Range("K105").Value = 10000
$ws.Range("M105").Value = -8253
$ws.Range("H115").Value = 34150
$ws.Range("J115").Value = 34150
$ws.Range("L115").Value = 34150
$ws.Range("N115").Value = -36500
$ws.Range("H122").Value = 1092.3334
$ws.Range("I122").Value = 1092.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3277.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -827.0001999999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3639.5386
$ws.Range("J126").Value = 7362.8
$ws.Range("L126").Value = 22088.4
$ws.Range("N126").Value = -27028.4
$ws.Range("H132").Value = 2205.9092
$ws.Range("I132").Value = 1783.7222
$ws.Range("J132").Value = 4105.75
$ws.Range("K132").Value = 5351.1666
$ws.Range("L132").Value = 12317.25
$ws.Range("M132").Value = -2821.1666
$ws.Range("N132").Value = -17377.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 50
$ws.Range("I24").Value = 50
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 150
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 80
$ws.Range("N24").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2998.6667
$ws.Range("I97").Value = 2998.6667
$ws.Range("K97").Value = 2998.6667
$ws.Range("M97").Value = -2502.6667
$ws.Range("H107").Value = 1225.5555
$ws.Range("I107").Value = 720.63635
$ws.Range("J107").Value = 1572.6875
$ws.Range("K107").Value = 720.63635
$ws.Range("L107").Value = 1572.6875
$ws.Range("M107").Value = 1199.36365
$ws.Range("N107").Value = -5412.6875
$ws.Range("H113").Value = 250000960
$ws.Range("I113").Value = 333334080
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 333334080
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = -333331910
$ws.Range("N113").Value = -5940
$ws.Range("H122").Value = 6380100
$ws.Range("I122").Value = 4322473.5
$ws.Range("K122").Value = 12967420.5
$ws.Range("M122").Value = -12964970.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3116.6667
$ws.Range("J7").Value = 5350
$ws.Range("L7").Value = 5350
$ws.Range("N7").Value = -5574
$ws.Range("H16").Value = 900.25
$ws.Range("I16").Value = 900.25
$ws.Range("K16").Value = 900.25
$ws.Range("M16").Value = -730.25
$ws.Range("H22").Value = 5051465
$ws.Range("I22").Value = 13889434
$ws.Range("J22").Value = 1196.9286
$ws.Range("K22").Value = 13889434
$ws.Range("L22").Value = 1196.9286
$ws.Range("M22").Value = -13889139
$ws.Range("N22").Value = -1786.9286
$ws.Range("H27").Value = 5051465
$ws.Range("I27").Value = 13889434
$ws.Range("J27").Value = 1196.9286
$ws.Range("K27").Value = 13889434
$ws.Range("L27").Value = 1196.9286
$ws.Range("M27").Value = -13889327
$ws.Range("N27").Value = -1410.9286
$ws.Range("H46").Value = 15874058
$ws.Range("I46").Value = 27778522
$ws.Range("J46").Value = 1438.6666
$ws.Range("K46").Value = 27778522
$ws.Range("L46").Value = 1438.6666
$ws.Range("M46").Value = -27778334
$ws.Range("N46").Value = -1814.6666
$ws.Range("H126").Value = 3116.6667
$ws.Range("J126").Value = 5350
$ws.Range("L126").Value = 16050
$ws.Range("N126").Value = -20990
$ws.Range("H132").Value = 21674156
$ws.Range("I132").Value = 33342750
$ws.Range("J132").Value = 3913.2856
$ws.Range("K132").Value = 100028250
$ws.Range("L132").Value = 11739.8568
$ws.Range("M132").Value = -100025720
$ws.Range("N132").Value = -16799.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 58824052
$ws.Range("J107").Value = 492
$ws.Range("L107").Value = 1476
$ws.Range("N107").Value = -5316
